# Commit: "Added more info about git diff."
#
# The paragraph under the "Git diff" heading currently ends with a run
# that explains what `git diff` does, followed by a lone trailing-space
# run. This adds a brand new run, right after that trailing space and
# still inside the same paragraph, with a further explanatory sentence -
# matching the plain (non-bold/underline) body-text formatting already
# used by the rest of the paragraph.

$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive existing text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This command helps us to compare the files*present working directory.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range

    # Pull the paragraph's own opening tag + <w:pPr> straight from the live
    # document so its identity (paraId, rsids, indentation, default run
    # formatting) is reused verbatim instead of being retyped by hand.
    $paraXml = $r.WordOpenXML
    $openMatch = [regex]::Match($paraXml, '<w:p\b[^>]*>(?:<w:pPr>.*?</w:pPr>)?')
    $paraOpen = $openMatch.Value

    $newSentence = "This command serves its purpose by showing us the changes made in the files in present working directory(which are not staged)."

    # Rebuild the paragraph: the existing "explanation" run, the existing
    # lone trailing-space run (both left byte-for-byte as they were), plus
    # the newly added run carrying the extra sentence with matching
    # (non-bold, sz 36 / 18pt, en-US) character formatting.
    $bodyXml = '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>This command helps us to compare the files/ documents in staging area (which are staged but not yet committed) vs the same files/documents in the present working directory.</w:t></w:r>' +
               '<w:r w:rsidR="00C56BCF"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
               '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr><w:t>' + $newSentence + '</w:t></w:r>'

    $newParaXml = $paraOpen + $bodyXml + '</w:p>'

    $xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
               '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # InsertXML replaces the contents of the range it's called on, so
    # targeting the whole paragraph range swaps in the rebuilt paragraph.
    $r.InsertXML($xmlFrag)
}
